# Insert a new weekly record for "Coliflor" (Macroferia Regional de Talca)
# at row 459, pushing the existing rows 459:478 down to 460:479.
# (New row 479 ends up holding what used to be row 478's data; the inserted
#  row 459 carries the newest observation.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(459).Insert()

$ws.Cells.Item(459, 1).Value  = 5
$ws.Cells.Item(459, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(459, 3).Value  = "Maule"
$ws.Cells.Item(459, 4).Value  = 45147
$ws.Cells.Item(459, 5).Value  = 7
$ws.Cells.Item(459, 6).Value  = 100112008
$ws.Cells.Item(459, 7).Value  = "Coliflor"
$ws.Cells.Item(459, 8).Value  = "Sin especificar"
$ws.Cells.Item(459, 9).Value  = "Primera"
$ws.Cells.Item(459, 10).Value = 3000
$ws.Cells.Item(459, 11).Value = 800
$ws.Cells.Item(459, 12).Value = 800
$ws.Cells.Item(459, 13).Value = 800
$ws.Cells.Item(459, 14).Value = "$/unidad"
$ws.Cells.Item(459, 15).Value = "Región del Maule"
$ws.Cells.Item(459, 16).Value = 800
$ws.Cells.Item(459, 17).Value = 1
$ws.Cells.Item(459, 18).Value = "Hortaliza"
